$d = $word.ActiveDocument
$nbsp = [char]160

function Replace-Text($findText, $replaceText) {
    $result = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        Write-Host "WARNING: not found ->" $findText
    }
    return $result
}

# 1. "Alle Kommentare werden auf englisch Verfasst " -> "Alle Kommentare werden auf Englisch verfasst "
Replace-Text "Alle Kommentare werden auf englisch Verfasst " "Alle Kommentare werden auf Englisch verfasst "

# 2. "Allgemein gilt: lieber ein Kommentar zu viel als eins zu wenig " -> "Allgemein gilt: lieber einen Kommentar zu viel als zu wenig "
Replace-Text "Allgemein gilt: lieber ein Kommentar zu viel als eins zu wenig " "Allgemein gilt: lieber einen Kommentar zu viel als zu wenig "

# 3. Remove lang=en-US from "Header " paragraph (pPr/rPr and run rPr)
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Header `r") {
        $p.Range.LanguageID = 0
        foreach ($r in $p.Range.Characters) { }
    }
}

# 4. funktionalität -> Funktionalität
Replace-Text "funktionalit" "Funktionalit"

# 5. "  */" with nbsp -> " *"
Replace-Text (" " + $nbsp + "*/") " *"

# 6. " /" + 70 stars -> " " + 71 stars  (slash -> extra star)
Replace-Text " /**********************************************************************" " ***********************************************************************"

# 7. Big sentence about "In einem Abschnitt..."
Replace-Text "In einem Abschnitt, also zB einer Funktion dem Loop dem setup etc werden die Kommentare die eine Zeile beschrieben alle auf einer linie untereinander verfasst" "In einem Abschnitt, also z.B. einer Funktion, dem Loop, dem Setup etc. werden die Kommentare, die eine Zeile beschrieben, alle auf einer Linie untereinander verfasst"

# 8. "Es wird alles ausgeschrieben außer..." -> add comma
Replace-Text "Es wird alles ausgeschrieben außer es gibt eine allgemein bekannte Abkürzung für das entsprechende Wort, dann ist dieses ebenfalls erlaubt " "Es wird alles ausgeschrieben, außer es gibt eine allgemein bekannte Abkürzung für das entsprechende Wort, dann ist dieses ebenfalls erlaubt "

# 9. "Es wird mit einem kleinen Buchstaben bekommen " -> "Es wird mit einem kleinen Buchstaben begonnen "
Replace-Text "Es wird mit einem kleinen Buchstaben bekommen " "Es wird mit einem kleinen Buchstaben begonnen "

# 10. "Varialennamen werden so gewählt das sie..." -> "Variablennamen werden so gewählt, dass sie..."
Replace-Text "Varialennamen werden so gewählt das sie möglichst eindeutig die Bedeutung der Variable beschreiben " "Variablennamen werden so gewählt, dass sie möglichst eindeutig die Bedeutung der Variable beschreiben "

# 11. "Jede Variable wird mit einem Kommentar beschrieben der Name..." -> add period + capitalize Der
Replace-Text "Jede Variable wird mit einem Kommentar beschrieben der Name könnte auch wenn es für den Programmierer eindeutig ist für Lesende nicht eindeutig sein " "Jede Variable wird mit einem Kommentar beschrieben. Der Name könnte auch wenn es für den Programmierer eindeutig ist für Lesende nicht eindeutig sein "

Write-Host "Done"
